$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.092.38"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "3.616.64"
$ws.Range("E3").Value = "  +3.34%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.77%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.652"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000305"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "4.191.83"
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.16%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.656.54"
$ws.Range("E16").Value = "  +4.53%  "
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "592.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "70.273.05"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.117"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").Value = "0.0₃0879"
$ws.Range("E35").Value = "  +8.67%  "
$ws.Range("D36").Value = "3.921.46"
$ws.Range("E36").Value = "  +5.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.59%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "522.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.393"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0456"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000252"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.47%  "
$ws.Range("E51").Value = "  +2.38%  "
